# dev-1.1.0 : export trace-list
# The xl_barcode template previously shipped with a hard-coded demo
# record (barcode value, part number, part model and part code values).
# Clear those sample values so the sheet is a clean, reusable template
# for exporting the trace-list - the field labels (PART NUMBER,
# PART MODEL, PART CODE, and the title) stay in place.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Barcode value (merged A2:E8) - was "*15011417A06C005*"
$ws.Range("A2").Value = ""

# PART CODE value (merged C9:E9) - was "15011417A06C005"
$ws.Range("C9").Value = ""

# PART NUMBER value (merged C10:E10) - was "TES"
$ws.Range("C10").Value = ""

# PART MODEL value (merged C11:E11) - was "TES"
$ws.Range("C11").Value = ""

# Trailing demo value (merged A12:E13) - was "OIL PAN"
$ws.Range("A12").Value = ""

# Move/restore the live selection to D18, matching the refreshed template
$ws.Range("D18").Select() | Out-Null
